$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on D-column cells whose new values would
# otherwise be auto-parsed as numbers (losing trailing zeros / going to
# scientific notation), so they stay literal strings like the source data.
$textCells = @("D5:D7", "D9:D15", "D17:D22", "D25:D29", "D31:D42", "D44:D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Column D (Price) updates ---
$ws.Range("D2").Value = '24.049.79'
$ws.Range("D3").Value = '1.668.44'
$ws.Range("D5").Value = '307.64'
$ws.Range("D6").Value = '0.9949'
$ws.Range("D7").Value = '0.3717'
$ws.Range("D9").Value = '47.80'
$ws.Range("D10").Value = '1.171'
$ws.Range("D11").Value = '0.07255'
$ws.Range("D12").Value = '0.9967'
$ws.Range("D13").Value = '20.50'
$ws.Range("D14").Value = '6.029'
$ws.Range("D15").Value = '6.731'
$ws.Range("D16").Value = '1.660.71'
$ws.Range("D17").Value = '0.00001099'
$ws.Range("D18").Value = '0.9948'
$ws.Range("D19").Value = '0.06707'
$ws.Range("D20").Value = '81.95'
$ws.Range("D21").Value = '16.40'
$ws.Range("D22").Value = '6.119'
$ws.Range("D24").Value = '24.000.25'
$ws.Range("D25").Value = '2.389'
$ws.Range("D26").Value = '3.382'
$ws.Range("D27").Value = '2.666'
$ws.Range("D28").Value = '152.12'
$ws.Range("D29").Value = '19.51'
$ws.Range("D30").Value = '1.842.66'
$ws.Range("D31").Value = '126.92'
$ws.Range("D32").Value = '6.385'
$ws.Range("D33").Value = '4.061'
$ws.Range("D34").Value = '0.9801'
$ws.Range("D35").Value = '1.726'
$ws.Range("D36").Value = '0.08383'
$ws.Range("D37").Value = '12.28'
$ws.Range("D38").Value = '8.947'
$ws.Range("D39").Value = '5.315'
$ws.Range("D40").Value = '0.06344'
$ws.Range("D41").Value = '1.288'
$ws.Range("D42").Value = '0.02318'
$ws.Range("D44").Value = '0.6094'
$ws.Range("D45").Value = '0.9948'
$ws.Range("D46").Value = '13.26'
$ws.Range("D47").Value = '3.809'
$ws.Range("D48").Value = '0.5937'
$ws.Range("D49").Value = '127.09'
$ws.Range("D50").Value = '2.001'
$ws.Range("D51").Value = '0.07095'

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = '  +16.15%  '
$ws.Range("E3").Value = '  +10.26%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("E5").Value = '  +10.16%  '
$ws.Range("E6").Value = '  +3.50%  '
$ws.Range("E7").Value = '  +4.18%  '
$ws.Range("E8").Value = '  +10.88%  '
$ws.Range("E9").Value = '  +20.72%  '
$ws.Range("E10").Value = '  +6.38%  '
$ws.Range("E11").Value = '  +7.98%  '
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("E13").Value = '  +10.57%  '
$ws.Range("E14").Value = '  +8.49%  '
$ws.Range("E15").Value = '  +7.69%  '
$ws.Range("E16").Value = '  +10.50%  '
$ws.Range("E17").Value = '  +7.07%  '
$ws.Range("E18").Value = '  +3.29%  '
$ws.Range("E19").Value = '  +11.06%  '
$ws.Range("E20").Value = '  +17.22%  '
$ws.Range("E21").Value = '  +11.03%  '
$ws.Range("E22").Value = '  +10.05%  '
$ws.Range("E23").Value = '  +6.38%  '
$ws.Range("E24").Value = '  +15.50%  '
$ws.Range("E25").Value = '  +3.90%  '
$ws.Range("E26").Value = '  -8.01%  '
$ws.Range("E27").Value = '  +24.20%  '
$ws.Range("E28").Value = '  +3.91%  '
$ws.Range("E29").Value = '  +11.81%  '
$ws.Range("E30").Value = '  +10.71%  '
$ws.Range("E31").Value = '  +9.64%  '
$ws.Range("E32").Value = '  +26.73%  '
$ws.Range("E33").Value = '  +2.12%  '
$ws.Range("E34").Value = '  +18.98%  '
$ws.Range("E35").Value = '  +18.49%  '
$ws.Range("E36").Value = '  +4.97%  '
$ws.Range("E37").Value = '  +17.28%  '
$ws.Range("E38").Value = '  +19.27%  '
$ws.Range("E39").Value = '  +10.56%  '
$ws.Range("E40").Value = '  +9.92%  '
$ws.Range("E42").Value = '  +12.88%  '
$ws.Range("E43").Value = '  +10.47%  '
$ws.Range("E44").Value = '  +15.16%  '
$ws.Range("E45").Value = '  +3.27%  '
$ws.Range("E46").Value = '  +7.72%  '
$ws.Range("E47").Value = '  +7.61%  '
$ws.Range("E48").Value = '  +12.73%  '
$ws.Range("E49").Value = '  +5.63%  '
$ws.Range("E50").Value = '  +7.95%  '
$ws.Range("E51").Value = '  +9.57%  '
